# aggiornamento fino a 6 gennaio 2022
# Append 27 new daily rows (465-491) to Sheet1, extending the data range
# from A1:D464 to A1:D491.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds dates formatted like the existing rows (style copied from
# the last existing data row, A464, which carries the date number format).
$ws.Range("A464").Copy()

# New data: date-serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$data = @(
    @(44539, 10,  26, 151.7362124306974),
    @(44540,  8,  33, 192.5882696235775),
    @(44541,  2,  30, 175.0802451123432),
    @(44542, 12,  39, 227.6043186460461),
    @(44543,  9,  46, 268.4563758389261),
    @(44544,  9,  53, 309.3084330318062),
    @(44545,  2,  52, 303.4724248613948),
    @(44546,  3,  45, 262.6203676685147),
    @(44547,  4,  41, 239.276334986869),
    @(44548, 12,  51, 297.6364166909834),
    @(44550, 13,  52, 303.4724248613948),
    @(44551,  3,  46, 268.4563758389261),
    @(44552,  3,  40, 233.4403268164575),
    @(44553, 14,  52, 303.4724248613948),
    @(44554,  9,  58, 338.4884738838634),
    @(44555, 26,  80, 466.880653632915),
    @(44556, 17,  85, 496.0606944849723),
    @(44557, 25,  97, 566.0927925299096),
    @(44558,  3,  97, 566.0927925299096),
    @(44559, 13, 107, 624.452874234024),
    @(44560, 27, 120, 700.3209804493727),
    @(44561, 78, 189, 1103.005544207762),
    @(44562, 36, 199, 1161.365625911876),
    @(44563, 18, 200, 1167.201634082288),
    @(44564, 10, 185, 1079.661511526116),
    @(44565,  2, 184, 1073.825503355705),
    @(44566, 16, 187, 1091.333527866939)
)

$startRow = 465
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    # Apply the date cell format (only) to the new A-column cell, matching
    # the style used by the rest of column A (s="2" / YYYY-MM-DD HH:MM:SS).
    $ws.Range("A$r").PasteSpecial(-4122)
}

Write-Host "Added rows $startRow to $($startRow + $data.Count - 1); new dimension should be A1:D$($startRow + $data.Count - 1)"
